$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# "English" appears multiple times in the document (hyperlinks + plain text).
# Only the single plain-text occurrence (paragraph 3, style P68B1DB1-Normal2)
# should be translated, so target it directly via the Paragraphs collection
# instead of a document-wide Find/Replace.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "English") {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = "ภาษาอังกฤษ"
}

Replace-All "Brief" "บทย่อ"

Replace-All "An email to partners in the the target country to invite them for a one-day seminar. It will be sent via customer.io" "An email to partners in the the target country to invite them for a one-day seminar. โดยมันจะถูกส่งผ่านทาง customer.io"

Replace-All "Target audience" "กลุ่มเป้าหมาย"

Replace-All "You’re invited to our Deriv Partner Seminar" "คุณได้รับเชิญให้เข้าร่วมงานสัมมนาหุ้นส่วน Deriv"

Replace-All "We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!" "เรารู้สึกตื่นเต้นที่จะแจ้งให้คุณทราบว่า ทีมพันธมิตร Deriv จะเยือน [CITY] ในเดือน [MONTH] เพื่อพบปะกับคุณผู้ซึ่งเป็นพันธมิตรที่มีค่าของเรา!"

Replace-All "In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. " "ในการสัมมนาหนึ่งวันนี้ เราจะให้การสนับสนุนด้านเทคนิคและการตลาด เปิดโอกาสให้คุณได้สร้างเครือข่ายกับพันธมิตรรายอื่นๆ ในระหว่างการรับประทานอาหารกลางวันแสนอร่อย รวมทั้งรับฟังความคิดเห็นของคุณเกี่ยวกับโครงการหุ้นส่วนพันธมิตรต่างๆ ของเรา นี่เป็นโอกาสของคุณที่จะแสดงความเห็นซึ่งจะช่วยให้เราวางแผนความพยายามสนับสนุนคุณในอนาคตให้ดียิ่งขึ้น "

Replace-All ". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!" ". Please note that attendance is confirmed on a first come, first served basis. เราหวังว่าจะได้พบเจอคุณที่นั่น!"

Replace-All "Send my details" "ส่งรายละเอียดของฉัน"
